# Append 5 new daily COVID-19 rows (2020-05-25 .. 2020-05-29) to the
# "Tabela1" table on the single worksheet, growing it from A1:J75 to A1:J80.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# New data rows: date-serial, Tested(all), Tested(daily), Positive(all),
# Positive(daily), Hospitalized, Intensive care, Discharged, Deaths(all),
# Deaths(daily)  -- columns A..J
$newRows = @(
    @(43976, 75770, 754, 1469, 0, 9, 2, 6, 108, 1),
    @(43977, 76579, 809, 1471, 2, 8, 2, 2, 108, 0),
    @(43978, 77210, 631, 1473, 2, 7, 2, 1, 108, 0),
    @(43979, 77916, 706, 1473, 0, 7, 2, 0, 108, 0),
    @(43980, 78529, 613, 1473, 0, 7, 2, 0, 108, 0)
)

$lastRow = $lo.Range().Row + $lo.Range().Rows.Count - 1

foreach ($rowValues in $newRows) {
    $newRowIndex = $lastRow + 1

    # Insert a formatted copy of the last data row right below it -- this
    # keeps the exact number formats / fonts / alignment used by the rest
    # of the table (matches the new cell styles Excel itself would reuse).
    $ws.Rows($lastRow).Copy()
    $ws.Rows($newRowIndex).Insert(-4121, 0)

    # Register the new row with the table so ref/autoFilter/dimension grow.
    $lo.ListRows.Add() | Out-Null

    # Overwrite the copied placeholder values with the real data for this row.
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($newRowIndex, $c).Value = $rowValues[$c - 1]
    }

    $lastRow = $newRowIndex
}

# Match the author's final selection: the newly added last row, A80:J80.
$ws.Range("A80:J80").Select()
